$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.505.53'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '1.871.66'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -2.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5081'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3902'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08355'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.85'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.215'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = '1.868.82'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.269'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06734'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -2.00%  '
$ws.Range("E22").Value = '  -1.43%  '
$ws.Range("D23").Value = '28.537.33'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.196'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.66%  '
$ws.Range("D26").Value = '2.082.46'
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.426'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1040'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.046'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.733'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.607'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06597'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2166'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.905'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.036'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.182'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.239'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("E44").Value = '  -1.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6008'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.678'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.009'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.213'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("E51").Value = '  -1.00%  '
